# Apply the "genotype" plate-map fill + re-activate that tab.
#
# Summary of the target edit (per the supplied diff):
#   - Fill the genotype plate map (B2:M9 on the "genotype" sheet) with the
#     new shared string "cep290_unknown" (the pipeline finally assigns the
#     96-well plate a genotype call instead of leaving it blank).
#   - Widen column D on that sheet a bit (new <cols> override) and move the
#     selection to B2.
#   - Make "genotype" the active/selected tab again (it had drifted to
#     "chem_perturbation").

$wb = $excel.ActiveWorkbook

$genotype = $wb.Worksheets.Item("genotype")

# Fill the whole plate-map body with the new call.
$genotype.Range("B2:M9").Value = "cep290_unknown"

# Give column D a bit more breathing room for the longer string.
$genotype.Columns.Item(4).ColumnWidth = 17

# Put the selection on B2 and bring this sheet to the front, which is what
# records it as the active tab for the workbook and clears tabSelected from
# whichever sheet had it before.
$genotype.Range("B2").Select()
$genotype.Activate()
